# Split the thesis-title run so the grammatical ending "ой" -> "ых"
# (параллельной -> параллельных) lives in its own run, matching the
# "split text into annot and main" commit: the single run holding
# "Исследование применимости сопрограмм в параллельной системах
# обработки данных" becomes three runs, with the corrected ending
# isolated in the middle one.

$d = $word.ActiveDocument

# Locate the whole title sentence first, so we don't depend on fixed
# character offsets elsewhere in the document.
$titleRange = $d.Content
$null = $titleRange.Find.Execute(
    "Исследование применимости сопрограмм в параллельной системах обработки данных",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentStart = $titleRange.Start
$sentEnd = $titleRange.End

# Within that sentence, find the "ой" ending right before " системах".
$scan = $d.Range($sentStart, $sentEnd)
$null = $scan.Find.Execute("ой системах", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
$endingStart = $scan.Start

# Narrow down to just the two-character ending "ой".
$ending = $d.Range($endingStart, $endingStart + 2)

# Fix the grammatical ending in place first (still one run at this point).
$ending.Text = "ых"

# Re-grab the (now 2-char, "ых") range and nudge a character-level
# property on and back off. That forces the run containing it to be
# split away from its neighbours, producing three runs: the prefix
# ("...параллельн"), the ending ("ых"), and the suffix (" системах
# обработки данных") — all keeping identical run formatting.
$ending = $d.Range($endingStart, $endingStart + 2)
$ending.Bold = 1
$ending.Bold = 0
